# Inventory.xlsx edit: add a "Threshold" column (quantity threshold used to
# flag low stock) between "StockCheckedon?" and "Stockless", and refresh a
# few Quantity values + the recomputed Stockless (YES/NO) flags that depend
# on the new Threshold column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at G. This shifts the existing
#    Stockless / Godown / Fromvendor(optional) columns from G/H/I to H/I/J,
#    carrying their data and header formatting with them.
$ws.Columns.Item(7).Insert()

# The freshly inserted column G inherited the date-style formatting from
# column F (its left neighbour). Reset the data cells (G2:G9) back to the
# default "Normal" style since Threshold holds plain numbers; the header
# cell G1 already correctly inherited the bold header style and is left
# untouched.
$ws.Range("G2:G9").Style = "Normal"

# 2. Header for the new column.
$ws.Cells.Item(1, 7).Value = "Threshold"

# 3. Threshold values (column G) for each product row.
$ws.Cells.Item(2, 7).Value = 10000
$ws.Cells.Item(3, 7).Value = 5000
$ws.Cells.Item(4, 7).Value = 14200
$ws.Cells.Item(5, 7).Value = 300
$ws.Cells.Item(6, 7).Value = 200
$ws.Cells.Item(7, 7).Value = 1000
$ws.Cells.Item(8, 7).Value = 1000
$ws.Cells.Item(9, 7).Value = 1000

# 4. Updated Quantity(KGS/MTS) values (column E) for a few rows.
$ws.Cells.Item(4, 5).Value = 10600
$ws.Cells.Item(5, 5).Value = 450
$ws.Cells.Item(9, 5).Value = 200

# 5. Stockless (now column H, after the insert) recomputed against the new
#    Threshold column: YES when Quantity < Threshold, otherwise NO.
$ws.Cells.Item(2, 8).Value = "YES"
$ws.Cells.Item(3, 8).Value = "NO"
$ws.Cells.Item(4, 8).Value = "YES"
$ws.Cells.Item(5, 8).Value = "NO"
$ws.Cells.Item(6, 8).Value = "NO"
$ws.Cells.Item(7, 8).Value = "NO"
$ws.Cells.Item(8, 8).Value = "NO"
$ws.Cells.Item(9, 8).Value = "YES"
